$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unmerge the old merged block A2:A5 ---
$ws.Range("A2:A5").UnMerge()

# --- Propagate the existing "A-column" style (currently only on A2:A5) down to A3:A21 ---
# Copying from A2 (which already carries the target style) preserves the existing
# style index instead of allocating brand-new ones.
$ws.Range("A2").Copy($ws.Range("A3:A21"))

# --- Fill in the B (answer) column first for rows 2-13, grouped by question, ---
# --- matching the original authoring order so shared-string indices line up. ---
$ws.Cells.Item(2, 2).Value = "Dùng eclipse"
$ws.Cells.Item(3, 2).Value = "Dùng VSCode"
$ws.Cells.Item(4, 2).Value = "Dùng Intellij"
$ws.Cells.Item(5, 2).Value = "Dùng website chính thứ của spring boot"

$ws.Cells.Item(6, 2).Value = "Chạy bằng STS"
$ws.Cells.Item(7, 2).Value = "Chạy bằng pring boot dashboard"
$ws.Cells.Item(8, 2).Value = "Answer Three"
$ws.Cells.Item(9, 2).Value = "Answer Four"

$ws.Cells.Item(10, 2).Value = "Không biết"
$ws.Cells.Item(11, 2).Value = "Không biết nốt"
$ws.Cells.Item(12, 2).Value = "foreach + CTRL + spacebar"
$ws.Cells.Item(13, 2).Value = "foreach +  spacebar"

# --- Now fill in the A (question) column for rows 2-13, one write per question group ---
$ws.Cells.Item(2, 1).Value = "Làm sao để tạo project Spring boot "
$ws.Cells.Item(3, 1).Value = "Làm sao để tạo project Spring boot "
$ws.Cells.Item(4, 1).Value = "Làm sao để tạo project Spring boot "
$ws.Cells.Item(5, 1).Value = "Làm sao để tạo project Spring boot "

$ws.Cells.Item(6, 1).Value = "Làm sao để chạy project Spring boot"
$ws.Cells.Item(7, 1).Value = "Làm sao để chạy project Spring boot"
$ws.Cells.Item(8, 1).Value = "Làm sao để chạy project Spring boot"
$ws.Cells.Item(9, 1).Value = "Làm sao để chạy project Spring boot"

$ws.Cells.Item(10, 1).Value = "Làm sao để chạy foreach project Spring boot"
$ws.Cells.Item(11, 1).Value = "Làm sao để chạy foreach project Spring boot"
$ws.Cells.Item(12, 1).Value = "Làm sao để chạy foreach project Spring boot"
$ws.Cells.Item(13, 1).Value = "Làm sao để chạy foreach project Spring boot"

# --- Rows 14-17 (endpoint question): question text first, then answers ---
$ws.Cells.Item(14, 1).Value = "Làm sao để tạo một endpoint project Spring boot"
$ws.Cells.Item(15, 1).Value = "Làm sao để tạo một endpoint project Spring boot"
$ws.Cells.Item(16, 1).Value = "Làm sao để tạo một endpoint project Spring boot"
$ws.Cells.Item(17, 1).Value = "Làm sao để tạo một endpoint project Spring boot"

$ws.Cells.Item(14, 2).Value = "Tạo trong controller"
$ws.Cells.Item(15, 2).Value = "Tạo service controller"
$ws.Cells.Item(16, 2).Value = "Tạo trong Model controller"
$ws.Cells.Item(17, 2).Value = "Tạo trong repository controller"

# --- Rows 18-20 (IDE question): question text first, then answers ---
$ws.Cells.Item(18, 1).Value = "IDE dùng để lập trình Spring boot"
$ws.Cells.Item(19, 1).Value = "IDE dùng để lập trình Spring boot"
$ws.Cells.Item(20, 1).Value = "IDE dùng để lập trình Spring boot"

$ws.Cells.Item(18, 2).Value = "VSCODE"
$ws.Cells.Item(19, 2).Value = "ECLIPSE"
$ws.Cells.Item(20, 2).Value = "Intelij"

# --- Correct-answer flags (column C) ---
$ws.Cells.Item(2, 3).Value = $true
$ws.Cells.Item(3, 3).Value = $true
$ws.Cells.Item(4, 3).Value = $true
$ws.Cells.Item(5, 3).Value = $true

$ws.Cells.Item(6, 3).Value = $true
$ws.Cells.Item(7, 3).Value = $true
$ws.Cells.Item(8, 3).Value = $false
$ws.Cells.Item(9, 3).Value = $false

$ws.Cells.Item(10, 3).Value = $false
$ws.Cells.Item(11, 3).Value = $false
$ws.Cells.Item(12, 3).Value = $false
$ws.Cells.Item(13, 3).Value = $true

$ws.Cells.Item(14, 3).Value = $true
$ws.Cells.Item(15, 3).Value = $false
$ws.Cells.Item(16, 3).Value = $false
$ws.Cells.Item(17, 3).Value = $false

$ws.Cells.Item(18, 3).Value = $true
$ws.Cells.Item(19, 3).Value = $true
$ws.Cells.Item(20, 3).Value = $true

# --- Row 21: style-only placeholder row, no content ---
$ws.Cells.Item(21, 1).Value = ""

# --- Alignment: drop the horizontal=center, keep vertical=center, for the whole A column block ---
$ws.Range("A2:A21").VerticalAlignment = -4108
$ws.Range("A2:A21").HorizontalAlignment = 1

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 42.877604166666664
$ws.Columns.Item(2).ColumnWidth = 35.736979166666664

# --- Selection ---
$ws.Range("F10").Select()
